$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "航天发展"
$ws.Range("B2").Value = "航天发展"
$ws.Range("C2").Value = "航天发展"

$ws.Range("A3").Value = "中国卫星"
$ws.Range("B3").Value = "中国卫星"
$ws.Range("C3").Value = "平潭发展"

$ws.Range("A4").Value = "航天电子"
$ws.Range("B4").Value = "平潭发展"
$ws.Range("C4").Value = "中国卫星"

$ws.Range("A5").Value = "天际股份"
$ws.Range("B5").Value = "海南发展"
$ws.Range("C5").Value = "再升科技"

$ws.Range("A6").Value = "神剑股份"
$ws.Range("B6").Value = "航天电子"
$ws.Range("C6").Value = "浙江世宝"

$ws.Range("A7").Value = "通宇通讯"
$ws.Range("B7").Value = "天际股份"
$ws.Range("C7").Value = "神剑股份"

$ws.Range("A8").Value = "锋龙股份"
$ws.Range("B8").Value = "神剑股份"
$ws.Range("C8").Value = "雪人集团"

$ws.Range("A9").Value = "平潭发展"
$ws.Range("B9").Value = "再升科技"
$ws.Range("C9").Value = "天际股份"

$ws.Range("A10").Value = "海南发展"
$ws.Range("B10").Value = "锋龙股份"
$ws.Range("C10").Value = "海南发展"

$ws.Range("A11").Value = "再升科技"
$ws.Range("B11").Value = "永鼎股份"
$ws.Range("C11").Value = "通宇通讯"

$ws.Range("A12").Value = "永鼎股份"
$ws.Range("B12").Value = "通宇通讯"
$ws.Range("C12").Value = "航天电子"

$ws.Range("A13").Value = "雪人集团"
$ws.Range("B13").Value = "多氟多"
$ws.Range("C13").Value = "东百集团"

$ws.Range("A14").Value = "浙江世宝"
$ws.Range("B14").Value = "东百集团"
$ws.Range("C14").Value = "永鼎股份"

$ws.Range("A15").Value = "西部材料"
$ws.Range("B15").Value = "雪人集团"
$ws.Range("C15").Value = "安通控股"

$ws.Range("A16").Value = "东百集团"
$ws.Range("B16").Value = "江西铜业"
$ws.Range("C16").Value = "胜通能源"

$ws.Range("A17").Value = "多氟多"
$ws.Range("B17").Value = "金风科技"
$ws.Range("C17").Value = "锋龙股份"

$ws.Range("A18").Value = "中超控股"
$ws.Range("B18").Value = "浙江世宝"
$ws.Range("C18").Value = "天普股份"

$ws.Range("A19").Value = "华丰科技"
$ws.Range("B19").Value = "中国卫通"
$ws.Range("C19").Value = "福龙马"

$ws.Range("A20").Value = "中国中免"
$ws.Range("B20").Value = "中国中免"
$ws.Range("C20").Value = "中超控股"

$ws.Range("A21").Value = "福龙马"
$ws.Range("B21").Value = "中超控股"
$ws.Range("C21").Value = "西部材料"
